$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 with the new week's data (20 Feb update -> date shifts from 09-Feb-22 to 15-Feb-22)
$ws.Range("A8").Value = 44607
$ws.Range("B8").Value = 0.05
$ws.Range("C8").Value = 0.15
$ws.Range("D8").Value = 0.01
$ws.Range("E8").Value = 0.45
$ws.Range("F8").Value = 0.11

# Update view/selection state to match the new focus area
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("D8").Select()

$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()

$wb.Save()
